$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2026-01-03 / 2026-01-04 NFL games (rows 258-273)
$row0 = @(46025, "Tampa Bay Buccaneers", "Carolina Panthers", -3, 16, 14, -1)
$row1 = @(46025, "San Francisco 49ers", "Seattle Seahawks", 2.5, 3, 13, -7.5)
$row2 = @(46026, "Jacksonville Jaguars", "Tennessee Titans", -12.5, 41, 7, 21.5)
$row3 = @(46026, "Houston Texans", "Indianapolis Colts", -9.5, 38, 30, -1.5)
$row4 = @(46026, "Cincinnati Bengals", "Cleveland Browns", -9.5, 18, 20, -11.5)
$row5 = @(46026, "New York Giants", "Dallas Cowboys", 3, 34, 17, 20)
$row6 = @(46026, "Atlanta Falcons", "New Orleans Saints", -3.5, 19, 17, -1.5)
$row7 = @(46026, "Minnesota Vikings", "Green Bay Packers", -13.5, 16, 3, -0.5)
$row8 = @(46026, "Los Angeles Rams", "Arizona Cardinals", -14.5, 37, 20, 2.5)
$row9 = @(46026, "Chicago Bears", "Detroit Lions", -3, 16, 19, -6)
$row10 = @(46026, "Philadelphia Eagles", "Washington Commanders", -3, 17, 24, -10)
$row11 = @(46026, "Las Vegas Raiders", "Kansas City Chiefs", 3.5, 14, 12, 5.5)
$row12 = @(46026, "Buffalo Bills", "New York Jets", -13.5, 35, 8, 13.5)
$row13 = @(46026, "New England Patriots", "Miami Dolphins", -14.5, 38, 10, 13.5)
$row14 = @(46026, "Denver Broncos", "Los Angeles Chargers", -15.5, 19, 3, 0.5)
$row15 = @(46026, "Pittsburgh Steelers", "Baltimore Ravens", 4.5, 26, 24, 6.5)

$data = @($row0, $row1, $row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15)

$startRow = 258
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}